# Apply the "add 2022-Q4 data" change:
#  1. Insert a brand-new worksheet "2022-Q4" right after "总计" (copied from
#     the "2022-Q3" sheet so that number formats / column widths / borders
#     match the other quarter sheets), then fill it with the 2022-Q4 numbers.
#  2. Insert a new row 2 in "总计" for the 2022-Q4 summary line, pushing the
#     existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

# Helper: write $val into $addr on worksheet $ws as a *text* cell, without
# introducing a new cell style (mirrors typing into a cell whose column is
# already formatted as General/Text in the source data).
function Set-TextCell($ws, $addr, [string]$val) {
    $escaped = $val -replace '"', '""'
    $c = $ws.Range($addr)
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Index = $q3Sheet.Index

$q3Sheet.Copy($q3Sheet, $null)
$newSheet = $wb.Worksheets.Item($q3Index)
$newSheet.Name = "2022-Q4"

Set-TextCell $newSheet "B2" "011685"
Set-TextCell $newSheet "C2" "创金合信先进装备股票A"
Set-TextCell $newSheet "D2" "0.24"
Set-TextCell $newSheet "E2" "80.29"
Set-TextCell $newSheet "F2" "7.23"
Set-TextCell $newSheet "G2" "0.0174"
$newSheet.Range("H2").Value = 7

Set-TextCell $newSheet "B3" "011686"
Set-TextCell $newSheet "C3" "创金合信先进装备股票C"
Set-TextCell $newSheet "D3" "0.18"
Set-TextCell $newSheet "E3" "80.29"
Set-TextCell $newSheet "F3" "7.23"
Set-TextCell $newSheet "G3" "0.0130"
$newSheet.Range("H3").Value = 7

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 summary row into "总计"
# ---------------------------------------------------------------------
$totalSheet.Range("A2:D2").Insert(-4121)   # xlShiftDown

$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A2").Value = 0
Set-TextCell $totalSheet "B2" "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.03

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
